$wb = $excel.ActiveWorkbook

# --- Update the panel query timestamps on the existing "data" sheet ---
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("F2").Value = "2021-10-05 14:22:58.466994"
$wsData.Range("F3").Value = "2021-10-05 14:22:58.467002"
$wsData.Range("F4").Value = "2021-10-05 14:22:58.467006"
$wsData.Range("F5").Value = "2021-10-05 14:22:58.467009"
$wsData.Range("F6").Value = "2021-10-05 14:22:58.467011"

# --- Add the new "metadata" sheet ---
$wsMeta = $wb.Worksheets.Add()
$wsMeta.Name = "metadata"

# Header row (row 1)
$wsMeta.Range("B1").Value = "data_name"
$wsMeta.Range("C1").Value = "data_id"
$wsMeta.Range("D1").Value = "data_version"
$wsMeta.Range("E1").Value = "data_version_created"
$wsMeta.Range("F1").Value = "panel_query_time"
$wsMeta.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$wsMeta.Range("A2").Value = 0
$wsMeta.Range("B2").Value = "Thyroid cancer pertinent cancer susceptibility"
$wsMeta.Range("C2").Value = 421

# D2 looks like a plain number, so force it to stay text (leading apostrophe, like Excel)
$wsMeta.Range("D2").Value = "'1.3"
$wsMeta.Range("E2").Value = "2021-09-24T11:00:38.835184Z"
$wsMeta.Range("F2").Value = "2021-10-05 14:22:58.463668"
$wsMeta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/421/?format=json"

# Formatting matching the bold/bordered/centered header style used on the "data" sheet
$headerRange = $wsMeta.Range("B1:G1")
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true

$idCell = $wsMeta.Range("A2")
$idCell.HorizontalAlignment = -4108
$idCell.VerticalAlignment = -4160
$idCell.Borders.LineStyle = 1
$idCell.Font.Bold = $true

# --- Move the new sheet so it comes right after "data" ---
$wsMeta.Move($null, $wb.Worksheets.Item("data"))
